$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- Add the extra (non-hidden, duplicate-with-suffix) AutoFilter defined name ---
# Mirrors the existing Blad1!_FilterDatabase name with a new "_0" suffixed name
# pointing at the same range, scoped to this worksheet.
$ws.Names.Add("_xlnm._FilterDatabase_0", "=Blad1!`$A`$1:`$J`$4")

# --- New data rows 5, 6 and 7 (BGT_OBW_opslagtank / BGT_OBW_overkapping / BGT_PND_pand) ---

# Row 5
$ws.Range("A5").Value = "BGT"
$ws.Range("B5").Value = "BGT_OBW_opslagtank"
$ws.Range("C5").Value = "gebouw"
$ws.Range("D5").Value = "vlak"
$ws.Range("E5").Formula = "=CONCATENATE(C5,""_"",D5,""<hoogteligging>"")"
$ws.Range("F5").Value = "lokaalid, bgttype, plustype"
$ws.Range("G5").Formula = "=C5"
$ws.Range("H5").Value = "gebouw"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 2
$ws.Rows.Item(5).RowHeight = 13.8

# Row 6
$ws.Range("A6").Value = "BGT"
$ws.Range("B6").Value = "BGT_OBW_overkapping"
$ws.Range("C6").Value = "gebouw"
$ws.Range("D6").Value = "vlak"
$ws.Range("E6").Formula = "=CONCATENATE(C6,""_"",D6,""<hoogteligging>"")"
$ws.Range("F6").Value = "lokaalid, bgttype, plustype"
$ws.Range("G6").Formula = "=C6"
$ws.Range("H6").Value = "gebouw"
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 4
$ws.Rows.Item(6).RowHeight = 13.8

# Row 7
$ws.Range("A7").Value = "BGT"
$ws.Range("B7").Value = "BGT_PND_pand"
$ws.Range("C7").Value = "gebouw"
$ws.Range("D7").Value = "vlak"
$ws.Range("E7").Formula = "=CONCATENATE(C7,""_"",D7,""<hoogteligging>"")"
$ws.Range("F7").Value = "lokaalid, bgttype, plustype"
$ws.Range("G7").Formula = "=C7"
$ws.Range("H7").Value = "gebouw"
$ws.Range("I7").Value = -1
$ws.Range("J7").Value = 3
$ws.Rows.Item(7).RowHeight = 15

# --- Move the active selection like the source document (F13 on Blad1) ---
[void]$ws.Range("F13").Select()
